$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 66
$ws.Range("C2").Value = 43

$ws.Range("B3").Value = 44
$ws.Range("C3").Value = 30

$ws.Range("B4").Value = 90
$ws.Range("C4").Value = 14

$ws.Range("B5").Value = 65
$ws.Range("C5").Value = 14

$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 88

$ws.Range("C7").Value = 61

$ws.Range("B8").Value = 74
$ws.Range("C8").Value = 69

$ws.Range("B9").Value = 41
$ws.Range("C9").Value = 29

$ws.Range("B10").Value = 18
$ws.Range("C10").Value = 31

$ws.Range("B11").Value = 14
$ws.Range("C11").Value = 33
